$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 46 (shifting nothing below it, since it's the last row)
# and have Excel carry the formatting down from the row above (45), the
# same way dragging the table's border down one row would.
$ws.Rows("46:46").Insert(-4121, 0)  # xlShiftDown, xlFormatFromLeftOrAbove

# A46 holds the same date label as A45 ("4.4.2020"). Force the cell to
# text first so Excel doesn't reinterpret the string as a date serial,
# then restore the (General) format copied from A45.
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value2 = "4.4.2020"
$ws.Range("A45").Copy() | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the rest of the new row's data. The Revision note (G) is typed
# in before the Task (E) is categorized, matching the order new entries
# were actually authored in the shared string table.
$ws.Range("B46").Value2 = 0.57291666666666663
$ws.Range("C46").Value2 = 0.61458333333333337
$ws.Range("D46").Formula = "=C46-B46"
$ws.Range("G46").Value2 = "Add packages for simulation"
$ws.Range("E46").Value2 = "Project"
$ws.Range("F46").Value2 = "Improve Code"

$ws.Range("E47").Select() | Out-Null
